$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('alicante','comunidad valenciana','españa'),
  @('andorra','andorra la vella','andorra'),
  @('avellaneda','buenos aires','argentina'),
  @('barquisimento','lara','venezuela'),
  @('barrio norte','buenos aires','argentina'),
  @('bogotá','cundinamarca','colombia'),
  @('buenos aires','buenos aires','argentina'),
  @('calafate','santa cruz','argentina'),
  @('campamento vespucio','salta','argentina'),
  @('capital federal','buenos aires','argentina'),
  @('caracas','distrito capital','venezuela'),
  @('cipolletti','río negro','argentina'),
  @('ciudad de méxico','ciudad de méxico','méxico'),
  @('comodoro rivadavia','chubut','argentina'),
  @('córdoba','córdoba','argentina'),
  @('coronel moldes','salta','argentina'),
  @('cumbrecita','córdoba','argentina'),
  @('curuzú cuatiá','corrientes','argentina'),
  @('embarcacion','salta','argentina'),
  @('ensenada','buenos aires','argentina'),
  @('florencia','caqueta','colombia'),
  @('granada','andalucía','españa'),
  @('guadalajara','jalisco','méxico'),
  @('humberto primo','santa fe','argentina'),
  @('jesus maria','córdoba','argentina'),
  @('la plata','buenos aires','argentina'),
  @('lisboa','lisboa','portugal'),
  @('lomas de zamora','buenos aires','argentina'),
  @('lund','escania','suecia'),
  @('málaga','andalucía','españa'),
  @('mar del plata','buenos aires','argentina'),
  @('maternidad sardá','buenos aires','argentina'),
  @('medellín','antioquia','colombia'),
  @('mendoza','mendoza','argentina'),
  @('merlo','buenos aires','argentina'),
  @('moisés ville','santa fe','argentina'),
  @('monte buey','cordoba','argentina'),
  @('monte grande','buenos aires','argentina'),
  @('montevideo','montevideo','uruguay'),
  @('moreno','buenos aires','argentina'),
  @('munro','buenos aires','argentina'),
  @('neuquén','neuquén','argentina'),
  @('palermo','buenos aires','argentina'),
  @('paraná','entre ríos','argentina'),
  @('parque patricios','buenos aires','argentina'),
  @('pasaje pujoi','corrientes','argentina'),
  @('pergamino','buenos aires','argentina'),
  @('perito moreno','santa cruz','argentina'),
  @('posadas','misiones','argentina'),
  @('rafaela','santa fe','argentina'),
  @('rio cuarto','cordoba','argentina'),
  @('rosario','santa fe','argentina'),
  @('san fernando','buenos aires','argentina'),
  @('san miguel','buenos aires','argentina'),
  @('san miguel de tucumán','tucumán','argentina'),
  @('san rafael','mendoza','argentina'),
  @('santa fe','santa fe','argentina'),
  @('sarmiento','chubut','argentina'),
  @('stroeder','buenos aires','argentina'),
  @('sunchales','santa fe','argentina'),
  @('temuco','araucanía','chile'),
  @('tierra del fuego','usuahia','argentina'),
  @('tostado','santa fe','argentina'),
  @('trelew','chubut','argentina'),
  @('valencia','comunidad valenciana','españa'),
  @('valladoi','castilla','españa'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
}

# Remove leftover sort state metadata (no active sort remains)
$ws.Sort.SortFields.Clear()

# Match final view/selection state
$ws.Range("B4").Select()
